# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Only the "K" column (G) values change for the affected rows; all other
# columns/values are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 2
    15 = 0
    16 = 1
    17 = 3
    18 = 1
    19 = 4
    20 = 3
    21 = 3
    22 = 2
    23 = 1
    24 = 2
    25 = 1
    28 = 1
    29 = 1
    30 = 3
    31 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
